$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 323, pushing the existing rows
# 323-355 down to 325-357 (new weekly data added at the top of the block).
$ws.Range("A323:A324").EntireRow.Insert()

# New row 323 ("Primera" quality) values for the new reporting date.
$ws.Range("A323").Value = 9
$ws.Range("B323").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C323").Value = "Metropolitana"
$ws.Range("D323").Value = 44918
$ws.Range("E323").Value = 13
$ws.Range("F323").Value = 100112017
$ws.Range("G323").Value = "Apio"
$ws.Range("H323").Value = "Americana (o)"
$ws.Range("I323").Value = "Primera"
$ws.Range("J323").Value = 70
$ws.Range("K323").Value = 9000
$ws.Range("L323").Value = 10000
$ws.Range("M323").Value = 9500
$ws.Range("N323").Value = "`$/docena de matas"
$ws.Range("O323").Value = "Región de Coquimbo"
$ws.Range("P323").Value = 1583
$ws.Range("Q323").Value = 6
$ws.Range("R323").Value = "Hortaliza"

# New row 324 ("Segunda" quality) values for the new reporting date.
$ws.Range("A324").Value = 9
$ws.Range("B324").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C324").Value = "Metropolitana"
$ws.Range("D324").Value = 44918
$ws.Range("E324").Value = 13
$ws.Range("F324").Value = 100112017
$ws.Range("G324").Value = "Apio"
$ws.Range("H324").Value = "Americana (o)"
$ws.Range("I324").Value = "Segunda"
$ws.Range("J324").Value = 52
$ws.Range("K324").Value = 7000
$ws.Range("L324").Value = 7000
$ws.Range("M324").Value = 7000
$ws.Range("N324").Value = "`$/docena de matas"
$ws.Range("O324").Value = "Región de Coquimbo"
$ws.Range("P324").Value = 1167
$ws.Range("Q324").Value = 6
$ws.Range("R324").Value = "Hortaliza"
